$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Update the "display.title" value for authorizations from "Authorization" to "Authorizations"
$ws.Range("C5").Value = "Authorizations"

# Make the "settings" sheet the active tab and select C5 (matches author's last view state)
$ws.Activate()
$ws.Range("C5").Select()
